$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 8.949600333333334
$ws.Cells.Item(2, 8).Value = 26.848801
$ws.Cells.Item(2, 9).Value = 0.2240220633121465
$ws.Cells.Item(2, 10).Value = 0.2240220633121465
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 26.31197333333334
$ws.Cells.Item(2, 14).Value = 78.93592000000001
$ws.Cells.Item(2, 15).Value = 0.2261559208386891
$ws.Cells.Item(2, 16).Value = 0.2261559208386891
$ws.Cells.Item(2, 17).Value = 235.4816453146578
$ws.Cells.Item(2, 18).Value = 2119.33480783192
$ws.Cells.Item(2, 19).Value = 0.0506639160165416
$ws.Cells.Item(2, 20).Value = 0.0506639160165416

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 8.949600333333334
$ws.Cells.Item(3, 8).Value = 26.848801
$ws.Cells.Item(3, 9).Value = 0.2240220633121465
$ws.Cells.Item(3, 10).Value = 0.2240220633121465
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 53.74150833333334
$ws.Cells.Item(3, 14).Value = 161.224525
$ws.Cells.Item(3, 15).Value = 0.461917475759518
$ws.Cells.Item(3, 16).Value = 0.461917475759518
$ws.Cells.Item(3, 17).Value = 480.9650208938362
$ws.Cells.Item(3, 18).Value = 4328.685188044526
$ws.Cells.Item(3, 19).Value = 0.1034797059995856
$ws.Cells.Item(3, 20).Value = 0.1034797059995856

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 8.949600333333334
$ws.Cells.Item(4, 8).Value = 26.848801
$ws.Cells.Item(4, 9).Value = 0.2240220633121465
$ws.Cells.Item(4, 10).Value = 0.2240220633121465
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 21.978693
$ws.Cells.Item(4, 14).Value = 65.93607899999999
$ws.Cells.Item(4, 15).Value = 0.1889106336220259
$ws.Cells.Item(4, 16).Value = 0.1889106336220259
$ws.Cells.Item(4, 17).Value = 196.700518199031
$ws.Cells.Item(4, 18).Value = 1770.304663791279
$ws.Cells.Item(4, 19).Value = 0.0423201499256112
$ws.Cells.Item(4, 20).Value = 0.04232014992561121

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 8.949600333333334
$ws.Cells.Item(5, 8).Value = 26.848801
$ws.Cells.Item(5, 9).Value = 0.2240220633121465
$ws.Cells.Item(5, 10).Value = 0.2240220633121465
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 14.31221833333333
$ws.Cells.Item(5, 14).Value = 42.936655
$ws.Cells.Item(5, 15).Value = 0.123015969779767
$ws.Cells.Item(5, 16).Value = 0.123015969779767
$ws.Cells.Item(5, 17).Value = 128.0886339667395
$ws.Cells.Item(5, 18).Value = 1152.797705700655
$ws.Cells.Item(5, 19).Value = 0.02755829137040806
$ws.Cells.Item(5, 20).Value = 0.02755829137040806

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 13.13839533333333
$ws.Cells.Item(6, 8).Value = 39.415186
$ws.Cells.Item(6, 9).Value = 0.328873952082703
$ws.Cells.Item(6, 10).Value = 0.3288739520827031
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 26.31197333333334
$ws.Cells.Item(6, 14).Value = 78.93592000000001
$ws.Cells.Item(6, 15).Value = 0.2261559208386891
$ws.Cells.Item(6, 16).Value = 0.2261559208386891
$ws.Cells.Item(6, 17).Value = 345.6971076534579
$ws.Cells.Item(6, 18).Value = 3111.27396888112
$ws.Cells.Item(6, 19).Value = 0.07437679147312262
$ws.Cells.Item(6, 20).Value = 0.07437679147312264

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 13.13839533333333
$ws.Cells.Item(7, 8).Value = 39.415186
$ws.Cells.Item(7, 9).Value = 0.328873952082703
$ws.Cells.Item(7, 10).Value = 0.3288739520827031
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 53.74150833333334
$ws.Cells.Item(7, 14).Value = 161.224525
$ws.Cells.Item(7, 15).Value = 0.461917475759518
$ws.Cells.Item(7, 16).Value = 0.461917475759518
$ws.Cells.Item(7, 17).Value = 706.0771822929612
$ws.Cells.Item(7, 18).Value = 6354.69464063665
$ws.Cells.Item(7, 19).Value = 0.1519126257890988
$ws.Cells.Item(7, 20).Value = 0.1519126257890989

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 13.13839533333333
$ws.Cells.Item(8, 8).Value = 39.415186
$ws.Cells.Item(8, 9).Value = 0.328873952082703
$ws.Cells.Item(8, 10).Value = 0.3288739520827031
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 21.978693
$ws.Cells.Item(8, 14).Value = 65.93607899999999
$ws.Cells.Item(8, 15).Value = 0.1889106336220259
$ws.Cells.Item(8, 16).Value = 0.1889106336220259
$ws.Cells.Item(8, 17).Value = 288.764757543966
$ws.Cells.Item(8, 18).Value = 2598.882817895694
$ws.Cells.Item(8, 19).Value = 0.06212778666972321
$ws.Cells.Item(8, 20).Value = 0.06212778666972323

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 13.13839533333333
$ws.Cells.Item(9, 8).Value = 39.415186
$ws.Cells.Item(9, 9).Value = 0.328873952082703
$ws.Cells.Item(9, 10).Value = 0.3288739520827031
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 14.31221833333333
$ws.Cells.Item(9, 14).Value = 42.936655
$ws.Cells.Item(9, 15).Value = 0.123015969779767
$ws.Cells.Item(9, 16).Value = 0.123015969779767
$ws.Cells.Item(9, 17).Value = 188.0395825603144
$ws.Cells.Item(9, 18).Value = 1692.35624304283
$ws.Cells.Item(9, 19).Value = 0.04045674815075833
$ws.Cells.Item(9, 20).Value = 0.04045674815075834

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 11.041444
$ws.Cells.Item(10, 8).Value = 33.124332
$ws.Cells.Item(10, 9).Value = 0.2763840813776585
$ws.Cells.Item(10, 10).Value = 0.2763840813776586
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 26.31197333333334
$ws.Cells.Item(10, 14).Value = 78.93592000000001
$ws.Cells.Item(10, 15).Value = 0.2261559208386891
$ws.Cells.Item(10, 16).Value = 0.2261559208386891
$ws.Cells.Item(10, 17).Value = 290.5221800894934
$ws.Cells.Item(10, 18).Value = 2614.699620805441
$ws.Cells.Item(10, 19).Value = 0.06250589642911955
$ws.Cells.Item(10, 20).Value = 0.06250589642911956

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 11.041444
$ws.Cells.Item(11, 8).Value = 33.124332
$ws.Cells.Item(11, 9).Value = 0.2763840813776585
$ws.Cells.Item(11, 10).Value = 0.2763840813776586
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 53.74150833333334
$ws.Cells.Item(11, 14).Value = 161.224525
$ws.Cells.Item(11, 15).Value = 0.461917475759518
$ws.Cells.Item(11, 16).Value = 0.461917475759518
$ws.Cells.Item(11, 17).Value = 593.3838547380334
$ws.Cells.Item(11, 18).Value = 5340.4546926423
$ws.Cells.Item(11, 19).Value = 0.1276666372100812
$ws.Cells.Item(11, 20).Value = 0.1276666372100813

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 11.041444
$ws.Cells.Item(12, 8).Value = 33.124332
$ws.Cells.Item(12, 9).Value = 0.2763840813776585
$ws.Cells.Item(12, 10).Value = 0.2763840813776586
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 21.978693
$ws.Cells.Item(12, 14).Value = 65.93607899999999
$ws.Cells.Item(12, 15).Value = 0.1889106336220259
$ws.Cells.Item(12, 16).Value = 0.1889106336220259
$ws.Cells.Item(12, 17).Value = 242.676507952692
$ws.Cells.Item(12, 18).Value = 2184.088571574228
$ws.Cells.Item(12, 19).Value = 0.05221189193609505
$ws.Cells.Item(12, 20).Value = 0.05221189193609506

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 11.041444
$ws.Cells.Item(13, 8).Value = 33.124332
$ws.Cells.Item(13, 9).Value = 0.2763840813776585
$ws.Cells.Item(13, 10).Value = 0.2763840813776586
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 14.31221833333333
$ws.Cells.Item(13, 14).Value = 42.936655
$ws.Cells.Item(13, 15).Value = 0.123015969779767
$ws.Cells.Item(13, 16).Value = 0.123015969779767
$ws.Cells.Item(13, 17).Value = 158.0275572432733
$ws.Cells.Item(13, 18).Value = 1422.24801518946
$ws.Cells.Item(13, 19).Value = 0.0339996558023627
$ws.Cells.Item(13, 20).Value = 0.0339996558023627

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 6.820198333333334
$ws.Cells.Item(14, 8).Value = 20.460595
$ws.Cells.Item(14, 9).Value = 0.1707199032274919
$ws.Cells.Item(14, 10).Value = 0.1707199032274919
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 26.31197333333334
$ws.Cells.Item(14, 14).Value = 78.93592000000001
$ws.Cells.Item(14, 15).Value = 0.2261559208386891
$ws.Cells.Item(14, 16).Value = 0.2261559208386891
$ws.Cells.Item(14, 17).Value = 179.4528766747111
$ws.Cells.Item(14, 18).Value = 1615.0758900724
$ws.Cells.Item(14, 19).Value = 0.03860931691990532
$ws.Cells.Item(14, 20).Value = 0.03860931691990532

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 6.820198333333334
$ws.Cells.Item(15, 8).Value = 20.460595
$ws.Cells.Item(15, 9).Value = 0.1707199032274919
$ws.Cells.Item(15, 10).Value = 0.1707199032274919
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 53.74150833333334
$ws.Cells.Item(15, 14).Value = 161.224525
$ws.Cells.Item(15, 15).Value = 0.461917475759518
$ws.Cells.Item(15, 16).Value = 0.461917475759518
$ws.Cells.Item(15, 17).Value = 366.5277455658195
$ws.Cells.Item(15, 18).Value = 3298.749710092375
$ws.Cells.Item(15, 19).Value = 0.07885850676075225
$ws.Cells.Item(15, 20).Value = 0.07885850676075225

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 6.820198333333334
$ws.Cells.Item(16, 8).Value = 20.460595
$ws.Cells.Item(16, 9).Value = 0.1707199032274919
$ws.Cells.Item(16, 10).Value = 0.1707199032274919
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 21.978693
$ws.Cells.Item(16, 14).Value = 65.93607899999999
$ws.Cells.Item(16, 15).Value = 0.1889106336220259
$ws.Cells.Item(16, 16).Value = 0.1889106336220259
$ws.Cells.Item(16, 17).Value = 149.899045367445
$ws.Cells.Item(16, 18).Value = 1349.091408307005
$ws.Cells.Item(16, 19).Value = 0.03225080509059644
$ws.Cells.Item(16, 20).Value = 0.03225080509059645

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 6.820198333333334
$ws.Cells.Item(17, 8).Value = 20.460595
$ws.Cells.Item(17, 9).Value = 0.1707199032274919
$ws.Cells.Item(17, 10).Value = 0.1707199032274919
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 14.31221833333333
$ws.Cells.Item(17, 14).Value = 42.936655
$ws.Cells.Item(17, 15).Value = 0.123015969779767
$ws.Cells.Item(17, 16).Value = 0.123015969779767
$ws.Cells.Item(17, 17).Value = 97.6121676233028
$ws.Cells.Item(17, 18).Value = 878.5095086097251
$ws.Cells.Item(17, 19).Value = 0.02100127445623789
$ws.Cells.Item(17, 20).Value = 0.02100127445623789
